$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 15.9605
$ws.Range("D3").Value = -8.605500000000001
$ws.Range("E3").Value = 16.51110000000001
$ws.Range("D4").Value = -7.382300000000002
$ws.Range("E9").Value = 17.4263
$ws.Range("C11").Value = -13.1597
$ws.Range("C12").Value = -10.7099
$ws.Range("D14").Value = -7.285399999999999
$ws.Range("C15").Value = -13.6824
$ws.Range("E15").Value = 16.4217
$ws.Range("E19").Value = 16.2076
$ws.Range("E20").Value = 16.1247
$ws.Range("E25").Value = 17.19820000000001
$ws.Range("D26").Value = -9.208799999999989
$ws.Range("C27").Value = -12.7304
$ws.Range("E27").Value = 16.68579999999999
$ws.Range("C28").Value = -13.96899999999999
$ws.Range("E28").Value = 16.08799999999999
$ws.Range("E30").Value = 15.0867
$ws.Range("C31").Value = -13.27740000000001
$ws.Range("D31").Value = -8.608600000000001
$ws.Range("C32").Value = -13.61940000000001
$ws.Range("E32").Value = 16.18389999999999
$ws.Range("D35").Value = -9.073999999999991
$ws.Range("C36").Value = -13.38110000000001
$ws.Range("D37").Value = -7.674399999999994
$ws.Range("C38").Value = -11.56650000000001
$ws.Range("D39").Value = -7.001600000000007
$ws.Range("D40").Value = -7.660699999999992
$ws.Range("E44").Value = 15.95510000000001
$ws.Range("D45").Value = -7.202100000000002
$ws.Range("C46").Value = -14.68959999999999
$ws.Range("E47").Value = 16.56950000000002
$ws.Range("D52").Value = -7.548599999999998
$ws.Range("C54").Value = -13.607
$ws.Range("C55").Value = -13.771
$ws.Range("C56").Value = -12.9081
$ws.Range("D57").Value = -8.623600000000003
$ws.Range("E58").Value = 16.9034
$ws.Range("E62").Value = 16.6361
$ws.Range("C67").Value = -10.6038
$ws.Range("C69").Value = -11.15009999999999
$ws.Range("C72").Value = -11.8236
$ws.Range("C73").Value = -12.9616
$ws.Range("E77").Value = 17.11720000000002
$ws.Range("E78").Value = 16.37820000000002
$ws.Range("D81").Value = -7.111699999999994
$ws.Range("C83").Value = -13.60440000000001
$ws.Range("D83").Value = -9.071699999999998
$ws.Range("E84").Value = 16.5481
$ws.Range("C86").Value = -14.02119999999999
$ws.Range("E89").Value = 17.2739
$ws.Range("C91").Value = -10.4314
$ws.Range("E91").Value = 17.83840000000001
$ws.Range("E92").Value = 18.12190000000001
$ws.Range("C93").Value = -11.369
$ws.Range("E96").Value = 15.8978
$ws.Range("C99").Value = -13.315
$ws.Range("D100").Value = -8.108300000000005
$ws.Range("D102").Value = -7.7447
$ws.Range("E102").Value = 16.5391
